$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the cached refresh-metadata blob stored in cell A1's comment
$commentText = "8RwAAB+LCAAAAAAAAAOlWVtvI0kV/istP8GD3W07s7mo0itfkmDhSxQ7ZLIvqN1diYu0u0xXdRK/LdpFO1oYCYndQcAuCCQEQmJ3JZYVuwOa/7IaZ4Yn/gKnLn2znRn3MIomXedWp06dOuerCnr7duob1zhkhAb7pWrFKhk4cKlHgsv9UsQvytW3Sm/b6ODWxf6xEzpTzEHYAK2A7d0ysl+acD7bM82bm5vKTb1Cw0uzZllV82GvO3QneOqUScC4E7i4lGh5r9cq2ajlTXuYO57DHaW5X+oMO5UWJm4baD0ncC5xWGlGjASYsYOAE04wE5ohdjhutXs/UAuza5W3KlVkrtBTyWZEfE/J5SQVXcvBtHhEptiuWdWdsrVbtuoja3uvZu1VrcpufeedWDERRF2H8SEOr4krCUPuTGdS3dq16tY2fFnIXCsEttIA2Gjgeyf4mjDstbDvs0IRMfUGNlwOqy4WTHAvo6sNvbkLR6Ezm4wI93ExN05Om8Y40L6kRmx0SEPsQvzeyKU+vhmEOqyjWRe4owkJ+bztzAvbOmU4HMxEkIqp2qhNA97wcchPZ7DX2INUAIbNwwgj8x5mqtQmzIVvEkTYsy8cn2WVckx0RsMrNnNc3IdzbAobN4FPHQ8SjhPGiZtOusJAxyGdgUWYvEl97xCsauE1jMRyJ4AQi2mblF6l3q1jIrmrcn9hT6cOj8VX6Gg4oTeDwJ8PozFzQzLGXrsZS6/lIXEgtXYrYpxOwYuUhBQtQ+n1zDn8gzO4zEFt7JKp4x/7EEdm18FQjoAaEacXhLeoH00DFru1REVnsKgRvk0WmYzRAPY3EHGnQSeI5VWk17LyCif0JplzlSHjkCE3mBvv+CpjWbgNtHgHVzlyU8QqD4kPLSK7HRlqPjGGE4z52qxQHCSq4aFoOnZz3o+mYzhhYzhm13JWhsyUjyBVId3BL9uCRlKWPyPL2pM/4EfCRgeBd79czEQwXWYuW5TpJRKCNflN3wmugHpG+KTfiNeyhoNUBO6VX+UhOLwz35lLchKlLA11AtePPKxqQie4kCkqfFObei8brZC6cMxt5ATz0XwGpZmRPQ4f+yVo1nuMhwAHSrZLo4CHc1E8kKlFX6fDonEgJ3D8jXUuQvzjCFDI/DAK3Bb1Np/NU9E5DQjf3EMahaoibq4ioyeKY8TaWJQZWfc31neLrImFhcSnAZ7SgLibRxuCLLz33mAhLD5VG2tgdb42lvehtavOJ876xmohQEjodYWmaTBGXSKTVR8PL6Nv3nNk2vjCiXyAbxy67GVSe5fJqMGulmWyJHQa+nEFtAU4ZoCOXW9acQE/CARYcelUEEwApWdDZGblBQhy8UFw2XWCywhgRlJXlulJ/RUtchQ6ARPLSVDFUileL4TiOqXQjq2K1yCSiaCKFwUuMpfk0AhPZzR0/B4EhhzqtNOQCdBIz+ETPYLe5mM3DrKZqiZaec9ix18nJpuUWoY48LpMLhGlkFiLwuGpTEpDYpU9OJZ+y/HJOFRVNW7l63iwYSk+jOuvWFxBrBjvAdzFoPt+H88FQE8Hmi5TthozVAKLQmoPT7Z2alu7O6Akx0iuuAU+UZ9IfGkcUbgOBlPQMwAZY0COe8acw1FQwdENqIBeVgMdxrVdTt+DFM0R8nyAGJcEOseKXMJI5e0eAM+JP89IqtV1qQuCi188Wnz91d1HzxbffHL30VeLn/7pv//6zeLrLxePPn35wd+fP/1QrU8Jo5Ez9rF0aNTc2bHqW5BaCQmJeJoSEHuRyyXt/Fzi4GSM9HVODloHndZRtylLSEKM1VUXMcVNcU6jdDhUi5ATyV00481XIvYoLkl6nONmupItbm7XOC+d5d+nqGLx4umfXzz9273aOmApvKru7tY3QV9w890qW7V70Fc3V/a18INyrZYRXpJBJ6rmJ3HqeHa9au1Wd7atalK+vSSD1wkts7SlkXNpLukpUkshoiQFsuOYKfN+hBlP2OokZAYqiHf/+NnLzz/OSenoakreCjgnkYuYzIwH0nT/ZGQMB6cnrQNjdDAUeZLyMnLK+CuE9ezJecomzeppy6XUOqoy90rF5RmPQhrNVipFSl0jubZmrHKXNGVsVl1LeWvk9a598midgl5JOwWkyRNGloZyHEXK8HXtevyT5988e/7s/Rf//tXdk88Xj/+4+Od7OTt6tuR+D6kLByQ7TDIZqpjuGksUdDaUMb2yfpjpEpoorkPHlASc2bUHW/IqpIcIdKvCnPyNOlPoXNKyDBvQlyjoew47uOX6sNp9ZOYJ4OjMgaZJ0ytkQlB1OQ3vfz793d1vv7x78sXLD/66+PAvi58/efH09y8/+4M6SXcff3H3+DNduZeLu/RFXEwVljPkS4driBNmiBZsfPvuL42AcgOQgxHJKvPtu7/OGBOOSoyRWgZkljiSd2FFNKss9IyMK4kPOb1ERfXxlmhL9URCNyY6I246yTtlYQrgpSEZ3+mMyhHDBgVQ9F1YSV44Vd5UT6uoNnm8bdWqNc1V3ogljB2WCf2RT8eOb8QM+Y6wJJLTerVCKivnO+oOmo1uKqKcGIQeDkUaqg8UI0PRJjosHsWplqEAF/CbG/ni6WdFbJWVWM6UM1O/oVw0PPGgvv7FISeBWlEYKpQT6Ff5YTQDTBs/tt3Plw+QGRjbV5AzC2zTcaed58M4w4XmlmcLguTLCqVZqlp1mHiuUai0L0KTDoGXe7SEcOiHd4WerglsiSkKz0EY0nBt9Uk5sVgPADFUFDONeCIj91SBZy/dq5gQV7zkQ13g9AppG/uYF3uVNlPtHqDfN9WFvS+q2mED39PBLHaDSMKSGsg+zYtE+X9f5lWyNcIQwJJ4xyv8lB7fP0/g2lrQG7UUqSgucjC7fus+JCHjD0Ul0F+Kcp5QzhXqfCjuTepDjs/t+gNFAAEza93MuRkfXa7+QEL9LpmSgrc7Kz7feSMQy9lMwbJOsUwRraWPbwE0ZixAURz/CNqGeg4pYk0lLNTSRF88QTJyOeFFHdseO9jDY6vsjnGtvOVZO+VdjOvlahX+d9xazbIeiAdMbRwqB8E3BScx4w1L/3Bp/w8RFcab8RwAAA=="
$ws.Range("A1").Comment.Text($commentText) | Out-Null

# 3. Update the number format for the data column (numFmtId 166: 0.000 -> ###0.000)
$ws.Range("B27:B36").NumberFormat = "###0.000"

# 4. Update label text in A11
$ws.Range("A11").Value = "Function Information"

# 5. Update the Skewness value in B20 (precision correction)
$ws.Range("B20").Value = -0.3785388757796309
